$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Set D4 (Make (m) / Buy (b) column) to "m" - it was missing the "number of camber shims" maker info
$ws.Range("D4").Value = "m"

# Update the active selection to E7 (as recorded by Excel when the file was last saved)
$ws.Range("E7").Select()
